$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Cells.Item(1, 1).Value = 'Rank'
$ws.Cells.Item(1, 2).Value = 'Game Name'
$ws.Cells.Item(1, 3).Value = 'Amount of Workshop Items as of 2024-04-22'

# Data rows
$ws.Cells.Item(2, 1).Value = 1
$ws.Cells.Item(2, 2).Value = 'Counter-Strike 2'
$ws.Cells.Item(2, 3).Value = 4648
$ws.Cells.Item(3, 1).Value = 2
$ws.Cells.Item(3, 2).Value = 'Dota 2'
$ws.Cells.Item(3, 3).Value = 32603
$ws.Cells.Item(4, 1).Value = 3
$ws.Cells.Item(4, 2).Value = 'Wallpaper Engine'
$ws.Cells.Item(4, 3).Value = 2187431
$ws.Cells.Item(5, 1).Value = 4
$ws.Cells.Item(5, 2).Value = 'Rust'
$ws.Cells.Item(5, 3).Value = 118623
$ws.Cells.Item(6, 1).Value = 5
$ws.Cells.Item(6, 2).Value = 'Team Fortress 2'
$ws.Cells.Item(6, 3).Value = 9527
$ws.Cells.Item(7, 1).Value = 6
$ws.Cells.Item(7, 2).Value = 'Warframe'
$ws.Cells.Item(7, 3).Value = 881
$ws.Cells.Item(8, 1).Value = 7
$ws.Cells.Item(8, 2).Value = 'Sid Meier’s Civilization® VI'
$ws.Cells.Item(8, 3).Value = 9316
$ws.Cells.Item(9, 1).Value = 8
$ws.Cells.Item(9, 2).Value = 'Unturned'
$ws.Cells.Item(9, 3).Value = 101670
$ws.Cells.Item(10, 1).Value = 9
$ws.Cells.Item(10, 2).Value = 'RimWorld'
$ws.Cells.Item(10, 3).Value = 35691
$ws.Cells.Item(11, 1).Value = 10
$ws.Cells.Item(11, 2).Value = 'Hearts of Iron IV'
$ws.Cells.Item(11, 3).Value = 48550
$ws.Cells.Item(12, 1).Value = 11
$ws.Cells.Item(12, 2).Value = 'Don''t Starve Together'
$ws.Cells.Item(12, 3).Value = 17179
$ws.Cells.Item(13, 1).Value = 12
$ws.Cells.Item(13, 2).Value = 'DayZ'
$ws.Cells.Item(13, 3).Value = 58704
$ws.Cells.Item(14, 1).Value = 13
$ws.Cells.Item(14, 2).Value = 'tModLoader'
$ws.Cells.Item(14, 3).Value = 6960
$ws.Cells.Item(15, 1).Value = 14
$ws.Cells.Item(15, 2).Value = 'Euro Truck Simulator 2'
$ws.Cells.Item(15, 3).Value = 22381
$ws.Cells.Item(16, 1).Value = 15
$ws.Cells.Item(16, 2).Value = 'Myth of Empires'
$ws.Cells.Item(16, 3).Value = 53
$ws.Cells.Item(17, 1).Value = 16
$ws.Cells.Item(17, 2).Value = 'ARK: Survival Evolved'
$ws.Cells.Item(17, 3).Value = 19625
$ws.Cells.Item(18, 1).Value = 17
$ws.Cells.Item(18, 2).Value = 'Squad'
$ws.Cells.Item(18, 3).Value = 463
$ws.Cells.Item(19, 1).Value = 18
$ws.Cells.Item(19, 2).Value = 'Mount & Blade II: Bannerlord'
$ws.Cells.Item(19, 3).Value = 619
$ws.Cells.Item(20, 1).Value = 19
$ws.Cells.Item(20, 2).Value = 'Slay the Spire'
$ws.Cells.Item(20, 3).Value = 947
$ws.Cells.Item(21, 1).Value = 20
$ws.Cells.Item(21, 2).Value = 'Left 4 Dead 2'
$ws.Cells.Item(21, 3).Value = 142550
$ws.Cells.Item(22, 1).Value = 21
$ws.Cells.Item(22, 2).Value = 'Project Zomboid'
$ws.Cells.Item(22, 3).Value = 26808
$ws.Cells.Item(23, 1).Value = 22
$ws.Cells.Item(23, 2).Value = 'Garry''s Mod'
$ws.Cells.Item(23, 3).Value = 1807968
$ws.Cells.Item(24, 1).Value = 23
$ws.Cells.Item(24, 2).Value = 'Crosshair X'
$ws.Cells.Item(24, 3).Value = 'N/A'
$ws.Cells.Item(25, 1).Value = 24
$ws.Cells.Item(25, 2).Value = 'VPet-Simulator'
$ws.Cells.Item(25, 3).Value = 601
$ws.Cells.Item(26, 1).Value = 25
$ws.Cells.Item(26, 2).Value = 'Europa Universalis IV'
$ws.Cells.Item(26, 3).Value = 12834
$ws.Cells.Item(27, 1).Value = 26
$ws.Cells.Item(27, 2).Value = 'Crusader Kings III'
$ws.Cells.Item(27, 3).Value = 8235
$ws.Cells.Item(28, 1).Value = 27
$ws.Cells.Item(28, 2).Value = 'Cities: Skylines'
$ws.Cells.Item(28, 3).Value = 340432
$ws.Cells.Item(29, 1).Value = 28
$ws.Cells.Item(29, 2).Value = 'Rocket League®'
$ws.Cells.Item(29, 3).Value = 608
$ws.Cells.Item(30, 1).Value = 29
$ws.Cells.Item(30, 2).Value = 'VTube Studio'
$ws.Cells.Item(30, 3).Value = 874
$ws.Cells.Item(31, 1).Value = 30
$ws.Cells.Item(31, 2).Value = 'Stellaris'
$ws.Cells.Item(31, 3).Value = 29413
$ws.Cells.Item(32, 1).Value = 31
$ws.Cells.Item(32, 2).Value = 'Cookie Clicker'
$ws.Cells.Item(32, 3).Value = 1245
$ws.Cells.Item(33, 1).Value = 32
$ws.Cells.Item(33, 2).Value = 'Arma 3'
$ws.Cells.Item(33, 3).Value = 144321
$ws.Cells.Item(34, 1).Value = 33
$ws.Cells.Item(34, 2).Value = 'Aimlabs'
$ws.Cells.Item(34, 3).Value = 42675
$ws.Cells.Item(35, 1).Value = 34
$ws.Cells.Item(35, 2).Value = 'Golf With Your Friends'
$ws.Cells.Item(35, 3).Value = 19108
$ws.Cells.Item(36, 1).Value = 35
$ws.Cells.Item(36, 2).Value = 'YoloMouse - Game Cursor Changer'
$ws.Cells.Item(36, 3).Value = 125
$ws.Cells.Item(37, 1).Value = 36
$ws.Cells.Item(37, 2).Value = 'Conan Exiles'
$ws.Cells.Item(37, 3).Value = 3902
$ws.Cells.Item(38, 1).Value = 37
$ws.Cells.Item(38, 2).Value = 'Brotato'
$ws.Cells.Item(38, 3).Value = 251
$ws.Cells.Item(39, 1).Value = 38
$ws.Cells.Item(39, 2).Value = 'Victoria 3'
$ws.Cells.Item(39, 3).Value = 4861
$ws.Cells.Item(40, 1).Value = 39
$ws.Cells.Item(40, 2).Value = '鬼谷八荒 Tale of Immortal'
$ws.Cells.Item(40, 3).Value = 9103
$ws.Cells.Item(41, 1).Value = 40
$ws.Cells.Item(41, 2).Value = 'Kenshi'
$ws.Cells.Item(41, 3).Value = 14736
$ws.Cells.Item(42, 1).Value = 41
$ws.Cells.Item(42, 2).Value = 'DSX'
$ws.Cells.Item(42, 3).Value = 2
$ws.Cells.Item(43, 1).Value = 42
$ws.Cells.Item(43, 2).Value = 'Mount & Blade: Warband'
$ws.Cells.Item(43, 3).Value = 401
$ws.Cells.Item(44, 1).Value = 43
$ws.Cells.Item(44, 2).Value = 'Noita'
$ws.Cells.Item(44, 3).Value = 1255
$ws.Cells.Item(45, 1).Value = 44
$ws.Cells.Item(45, 2).Value = 'American Truck Simulator'
$ws.Cells.Item(45, 3).Value = 10727
$ws.Cells.Item(46, 1).Value = 45
$ws.Cells.Item(46, 2).Value = 'Dying Light'
$ws.Cells.Item(46, 3).Value = 806
$ws.Cells.Item(47, 1).Value = 46
$ws.Cells.Item(47, 2).Value = 'Company of Heroes 2'
$ws.Cells.Item(47, 3).Value = 14909
$ws.Cells.Item(48, 1).Value = 47
$ws.Cells.Item(48, 2).Value = 'Planet Zoo'
$ws.Cells.Item(48, 3).Value = 104266
$ws.Cells.Item(49, 1).Value = 48
$ws.Cells.Item(49, 2).Value = 'People Playground'
$ws.Cells.Item(49, 3).Value = 709113
$ws.Cells.Item(50, 1).Value = 49
$ws.Cells.Item(50, 2).Value = 'SAO Utils 2: Progressive'
$ws.Cells.Item(50, 3).Value = 'N/A'
$ws.Cells.Item(51, 1).Value = 50
$ws.Cells.Item(51, 2).Value = 'Divinity: Original Sin 2 - Definitive Edition'
$ws.Cells.Item(51, 3).Value = 4704
$ws.Cells.Item(52, 1).Value = 51
$ws.Cells.Item(52, 2).Value = 'XCOM® 2'
$ws.Cells.Item(52, 3).Value = 8543
$ws.Cells.Item(53, 1).Value = 52
$ws.Cells.Item(53, 2).Value = 'They Are Billions'
$ws.Cells.Item(53, 3).Value = 4302
$ws.Cells.Item(54, 1).Value = 53
$ws.Cells.Item(54, 2).Value = 'Space Engineers'
$ws.Cells.Item(54, 3).Value = 552484
$ws.Cells.Item(55, 1).Value = 54
$ws.Cells.Item(55, 2).Value = 'Halo: The Master Chief Collection'
$ws.Cells.Item(55, 3).Value = 1143
$ws.Cells.Item(56, 1).Value = 55
$ws.Cells.Item(56, 2).Value = 'Dead Cells'
$ws.Cells.Item(56, 3).Value = 724
$ws.Cells.Item(57, 1).Value = 56
$ws.Cells.Item(57, 2).Value = 'F1® 23'
$ws.Cells.Item(57, 3).Value = 1284
$ws.Cells.Item(58, 1).Value = 57
$ws.Cells.Item(58, 2).Value = 'Age of Empires II (Retired)'
$ws.Cells.Item(58, 3).Value = 17415
$ws.Cells.Item(59, 1).Value = 58
$ws.Cells.Item(59, 2).Value = 'Tabletop Simulator'
$ws.Cells.Item(59, 3).Value = 82664
$ws.Cells.Item(60, 1).Value = 59
$ws.Cells.Item(60, 2).Value = 'Hero''s Adventure: Road to Passion'
$ws.Cells.Item(60, 3).Value = 771
$ws.Cells.Item(61, 1).Value = 60
$ws.Cells.Item(61, 2).Value = 'Transport Fever 2'
$ws.Cells.Item(61, 3).Value = 13049
$ws.Cells.Item(62, 1).Value = 61
$ws.Cells.Item(62, 2).Value = 'Call of Duty®: Black Ops III'
$ws.Cells.Item(62, 3).Value = 5255
$ws.Cells.Item(63, 1).Value = 62
$ws.Cells.Item(63, 2).Value = 'Farthest Frontier'
$ws.Cells.Item(63, 3).Value = 'N/A'
$ws.Cells.Item(64, 1).Value = 63
$ws.Cells.Item(64, 2).Value = 'Warhammer 40,000: Rogue Trader'
$ws.Cells.Item(64, 3).Value = 14
$ws.Cells.Item(65, 1).Value = 64
$ws.Cells.Item(65, 2).Value = 'Killing Floor 2'
$ws.Cells.Item(65, 3).Value = 3012
$ws.Cells.Item(66, 1).Value = 65
$ws.Cells.Item(66, 2).Value = 'A Dance of Fire and Ice'
$ws.Cells.Item(66, 3).Value = 18870
$ws.Cells.Item(67, 1).Value = 66
$ws.Cells.Item(67, 2).Value = 'MyDockFinder'
$ws.Cells.Item(67, 3).Value = 3660
$ws.Cells.Item(68, 1).Value = 67
$ws.Cells.Item(68, 2).Value = 'X4: Foundations'
$ws.Cells.Item(68, 3).Value = 829
$ws.Cells.Item(69, 1).Value = 68
$ws.Cells.Item(69, 2).Value = 'Total War: WARHAMMER II'
$ws.Cells.Item(69, 3).Value = 12719
$ws.Cells.Item(70, 1).Value = 69
$ws.Cells.Item(70, 2).Value = 'SAO Utils: Beta'
$ws.Cells.Item(70, 3).Value = 263
$ws.Cells.Item(71, 1).Value = 70
$ws.Cells.Item(71, 2).Value = 'Pummel Party'
$ws.Cells.Item(71, 3).Value = 1353
$ws.Cells.Item(72, 1).Value = 71
$ws.Cells.Item(72, 2).Value = '觅长生'
$ws.Cells.Item(72, 3).Value = 1512
$ws.Cells.Item(73, 1).Value = 72
$ws.Cells.Item(73, 2).Value = 'Company of Heroes 3'
$ws.Cells.Item(73, 3).Value = 509
$ws.Cells.Item(74, 1).Value = 73
$ws.Cells.Item(74, 2).Value = 'Human Fall Flat'
$ws.Cells.Item(74, 3).Value = 504748
$ws.Cells.Item(75, 1).Value = 74
$ws.Cells.Item(75, 2).Value = 'Call to Arms - Gates of Hell: Ostfront'
$ws.Cells.Item(75, 3).Value = 1494
$ws.Cells.Item(76, 1).Value = 75
$ws.Cells.Item(76, 2).Value = 'KovaaK''s'
$ws.Cells.Item(76, 3).Value = 32435
$ws.Cells.Item(77, 1).Value = 76
$ws.Cells.Item(77, 2).Value = 'Banana Shooter'
$ws.Cells.Item(77, 3).Value = 968
$ws.Cells.Item(78, 1).Value = 77
$ws.Cells.Item(78, 2).Value = 'Don''t Starve'
$ws.Cells.Item(78, 3).Value = 3200
$ws.Cells.Item(79, 1).Value = 78
$ws.Cells.Item(79, 2).Value = 'CarX Drift Racing Online'
$ws.Cells.Item(79, 3).Value = 870
$ws.Cells.Item(80, 1).Value = 79
$ws.Cells.Item(80, 2).Value = 'Teardown'
$ws.Cells.Item(80, 3).Value = 6652
$ws.Cells.Item(81, 1).Value = 80
$ws.Cells.Item(81, 2).Value = 'Football Manager 2020'
$ws.Cells.Item(81, 3).Value = 20132
$ws.Cells.Item(82, 1).Value = 81
$ws.Cells.Item(82, 2).Value = 'Barotrauma'
$ws.Cells.Item(82, 3).Value = 54353
$ws.Cells.Item(83, 1).Value = 82
$ws.Cells.Item(83, 2).Value = 'Stranded: Alien Dawn'
$ws.Cells.Item(83, 3).Value = 786
$ws.Cells.Item(84, 1).Value = 83
$ws.Cells.Item(84, 2).Value = 'The Elder Scrolls V: Skyrim'
$ws.Cells.Item(84, 3).Value = 27731
$ws.Cells.Item(85, 1).Value = 84
$ws.Cells.Item(85, 2).Value = 'Songs of Syx'
$ws.Cells.Item(85, 3).Value = 236
$ws.Cells.Item(86, 1).Value = 85
$ws.Cells.Item(86, 2).Value = 'House Flipper'
$ws.Cells.Item(86, 3).Value = 30089
$ws.Cells.Item(87, 1).Value = 86
$ws.Cells.Item(87, 2).Value = 'Workers & Resources: Soviet Republic'
$ws.Cells.Item(87, 3).Value = 9163
$ws.Cells.Item(88, 1).Value = 87
$ws.Cells.Item(88, 2).Value = 'Hydroneer'
$ws.Cells.Item(88, 3).Value = 57
$ws.Cells.Item(89, 1).Value = 88
$ws.Cells.Item(89, 2).Value = 'Portal 2'
$ws.Cells.Item(89, 3).Value = 948465
$ws.Cells.Item(90, 1).Value = 89
$ws.Cells.Item(90, 2).Value = 'Fisher Online'
$ws.Cells.Item(90, 3).Value = 578
$ws.Cells.Item(91, 1).Value = 90
$ws.Cells.Item(91, 2).Value = 'PlateUp!'
$ws.Cells.Item(91, 3).Value = 261
$ws.Cells.Item(92, 1).Value = 91
$ws.Cells.Item(92, 2).Value = 'Trove'
$ws.Cells.Item(92, 3).Value = 2188
$ws.Cells.Item(93, 1).Value = 92
$ws.Cells.Item(93, 2).Value = 'Age of Mythology: Extended Edition'
$ws.Cells.Item(93, 3).Value = 2223
$ws.Cells.Item(94, 1).Value = 93
$ws.Cells.Item(94, 2).Value = 'Wobbly Life'
$ws.Cells.Item(94, 3).Value = 59
$ws.Cells.Item(95, 1).Value = 94
$ws.Cells.Item(95, 2).Value = 'Kerbal Space Program'
$ws.Cells.Item(95, 3).Value = 105595
$ws.Cells.Item(96, 1).Value = 95
$ws.Cells.Item(96, 2).Value = 'Planet Coaster'
$ws.Cells.Item(96, 3).Value = 408730
$ws.Cells.Item(97, 1).Value = 96
$ws.Cells.Item(97, 2).Value = 'Library Of Ruina'
$ws.Cells.Item(97, 3).Value = 5920
$ws.Cells.Item(98, 1).Value = 97
$ws.Cells.Item(98, 2).Value = 'Scrap Mechanic'
$ws.Cells.Item(98, 3).Value = 480190
$ws.Cells.Item(99, 1).Value = 98
$ws.Cells.Item(99, 2).Value = 'Stormworks: Build and Rescue'
$ws.Cells.Item(99, 3).Value = 246550
$ws.Cells.Item(100, 1).Value = 99
$ws.Cells.Item(100, 2).Value = 'Age of Wonders 4'
$ws.Cells.Item(100, 3).Value = 719
$ws.Cells.Item(101, 1).Value = 100
$ws.Cells.Item(101, 2).Value = 'Company of Heroes'
$ws.Cells.Item(101, 3).Value = 3337
